$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (ECs -> Csf3r target MuSCs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06796566666666666
$ws.Range("H2").Value = 0.203897
$ws.Range("I2").Value = 0.5354438025210083
$ws.Range("J2").Value = 0.5354438025210083
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.031135
$ws.Range("N2").Value = 0.093405
$ws.Range("Q2").Value = 0.002116111031666667
$ws.Range("R2").Value = 0.019044999285
$ws.Range("S2").Value = 0.5354438025210083
$ws.Range("T2").Value = 0.5354438025210083

# Row 3 updates
$ws.Range("G3").Value = 0.05896766666666667
$ws.Range("H3").Value = 0.176903
$ws.Range("I3").Value = 0.4645561974789916
$ws.Range("J3").Value = 0.4645561974789916
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.031135
$ws.Range("N3").Value = 0.093405
$ws.Range("Q3").Value = 0.001835958301666667
$ws.Range("R3").Value = 0.016523624715
$ws.Range("S3").Value = 0.4645561974789916
$ws.Range("T3").Value = 0.4645561974789916
